$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-obsolete "ability after guard activated" entries in A6:A7
# without shifting any other rows (row 13 keeps its row number).
$ws.Range("A6:A7").ClearContents()

# Move the active selection to where row 7 used to be (now row 6).
$ws.Range("A6").Select()
